# The document has a title-page layout: each section has a "first page"
# header/footer pair (header1.xml / footer1.xml) plus a "default" (primary)
# header/footer pair (header2.xml / footer2.xml). Both pairs contain an
# inline picture - the Pearson logo (PNG) in the footers, and the BTec logo
# (JPG) in the headers. The picture's internal <wp:docPr>/name attribute was
# swapped with its sibling's, i.e. the Pearson logo is now mis-labelled
# "image1.png" (it used to be "image2.png") and the BTec logo is now
# mis-labelled "image2.jpg" (it used to be "image1.jpg").

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# --- Footers (Pearson Edexcel logo, PNG): image2.png -> image1.png ---
$ftrFirst = $sec.Footers.Item(2)          # wdHeaderFooterFirstPage -> footer1.xml
$ftrFirst.Range.InlineShapes.Item(1).Name = "image1.png"

$ftrPrimary = $sec.Footers.Item(1)        # wdHeaderFooterPrimary   -> footer2.xml
$ftrPrimary.Range.InlineShapes.Item(1).Name = "image1.png"

# --- Headers (BTec logo, JPG): image1.jpg -> image2.jpg ---
$hdrFirst = $sec.Headers.Item(2)          # wdHeaderFooterFirstPage -> header1.xml
$hdrFirst.Range.InlineShapes.Item(1).Name = "image2.jpg"

$hdrPrimary = $sec.Headers.Item(1)        # wdHeaderFooterPrimary   -> header2.xml
$hdrPrimary.Range.InlineShapes.Item(1).Name = "image2.jpg"
